# Generate Report for Handoff
#
# Regenerating the localization-status handoff report bumped the "Latest HO
# Xliff Generate Date" / "Latest Handoff Datetime" timestamps for the batch
# of files that were "Ready for handoff", and flipped their handoff
# Priority marker from "ht" (human translation) to "mt" (machine
# translation) now that the handoff xliffs were (re)generated.

$wb = $excel.ActiveWorkbook

# Rows on every per-locale sheet (Overview / zh-cn / de-de) that belong to
# files with Status = "Ready for handoff" and Priority = "ht" prior to this
# run. Rows 12 and 15 are excluded: they have no Priority set (blank), so
# they are not part of this handoff batch.
$rows = @(7, 8, 9, 10, 11, 13, 14, 16)

# Overview sheet: "Latest HO Xliff Generate Date" column G.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-09-03 16:27:33"
}

# zh-cn sheet: Priority column E flips ht -> mt, and "Latest Handoff
# Datetime" column H gets the new zh-cn-specific generation timestamp.
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "mt"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-03 16:27:28"
}

# de-de sheet: Priority column E flips ht -> mt, and "Latest Handoff
# Datetime" column H gets the new de-de-specific generation timestamp
# (matches the Overview sheet's generate date for this batch).
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "mt"
    $dede.Cells.Item($r, 8).Value = "2016-09-03 16:27:33"
}
